$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# home_score (column D) and away_score (column E) values for rows 2-15,
# matching the home_team/away_team rows already present in the sheet.
$scores = @(
    @(2, 23.0, 26.0),
    @(3, 17.0, 21.0),
    @(4, 10.0, 44.0),
    @(5, 27.0, 24.0),
    @(6, 22.0, 37.0),
    @(7, 17.0, 21.0),
    @(8, 40.0, 6.0),
    @(9, 26.0, 14.0),
    @(10, 21.0, 22.0),
    @(11, 35.0, 38.0),
    @(12, 24.0, 37.0),
    @(13, 10.0, 27.0),
    @(14, 20.0, 23.0),
    @(15, 31.0, 28.0)
)

# Capture the formatting already used by the existing D2 score cell so the
# newly written score cells look consistent with it.
$scoreStyle = $ws.Cells.Item(2, 4).Style

foreach ($row in $scores) {
    $r = $row[0]
    $homeScore = $row[1]
    $awayScore = $row[2]
    $ws.Cells.Item($r, 4).Value = $homeScore
    $ws.Cells.Item($r, 5).Value = $awayScore
}

# Apply consistent formatting to every newly written score cell (D3:E15 plus
# the new E2 cell; D2 already had this formatting).
$ws.Range("E2").Style = $scoreStyle
$ws.Range("D3:E15").Style = $scoreStyle
